$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "98.499.59"
$ws.Cells.Item(2, 5).Value = "  -0.44%  "
$ws.Cells.Item(3, 4).Value = "3.380.11"
$ws.Cells.Item(3, 5).Value = "  +0.11%  "
$ws.Cells.Item(4, 5).Value = "  -0.06%  "
$ws.Cells.Item(5, 4).Value = "'258.87"
$ws.Cells.Item(5, 5).Value = "  -0.41%  "
$ws.Cells.Item(6, 4).Value = "'669.54"
$ws.Cells.Item(6, 5).Value = "  +6.43%  "
$ws.Cells.Item(7, 5).Value = "  +12.83%  "
$ws.Cells.Item(8, 4).Value = "'0.459"
$ws.Cells.Item(8, 5).Value = "  +16.92%  "
$ws.Cells.Item(9, 4).Value = "'1.11"
$ws.Cells.Item(9, 5).Value = "  +28.19%  "
$ws.Cells.Item(10, 5).Value = "  -0.01%  "
$ws.Cells.Item(11, 4).Value = "3.376.80"
$ws.Cells.Item(11, 5).Value = "  +0.09%  "
$ws.Cells.Item(12, 5).Value = "  +5.48%  "
$ws.Cells.Item(13, 4).Value = "'42.57"
$ws.Cells.Item(13, 5).Value = "  +17.90%  "
$ws.Cells.Item(14, 5).Value = "  +8.01%  "
$ws.Cells.Item(15, 4).Value = "98.291.72"
$ws.Cells.Item(15, 5).Value = "  -0.43%  "
$ws.Cells.Item(16, 2).Value = "Toncoin"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(16, 4).Value = "'5.63"
$ws.Cells.Item(16, 5).Value = "  +2.42%  "
$ws.Cells.Item(17, 2).Value = "WrappedEther"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(17, 4).Value = "3.373.77"
$ws.Cells.Item(17, 5).Value = "  -0.28%  "
$ws.Cells.Item(18, 2).Value = "Polkadot"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(18, 4).Value = "'7.63"
$ws.Cells.Item(18, 5).Value = "  +24.42%  "
$ws.Cells.Item(19, 2).Value = "Chainlink"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(19, 4).Value = "'17.02"
$ws.Cells.Item(19, 5).Value = "  +11.91%  "
$ws.Cells.Item(20, 2).Value = "SuiNetwork"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(20, 4).Value = "'3.60"
$ws.Cells.Item(20, 5).Value = "  +1.18%  "
$ws.Cells.Item(21, 2).Value = "BitcoinCash"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(21, 4).Value = "'530.61"
$ws.Cells.Item(21, 5).Value = "  +8.20%  "
$ws.Cells.Item(22, 2).Value = "Uniswap"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(22, 4).Value = "'10.54"
$ws.Cells.Item(22, 5).Value = "  +12.19%  "
$ws.Cells.Item(23, 2).Value = "Stellar"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(23, 4).Value = "'0.450"
$ws.Cells.Item(23, 5).Value = "  +60.53%  "
$ws.Cells.Item(24, 4).Value = "'0.0000214"
$ws.Cells.Item(24, 5).Value = "  +2.24%  "
$ws.Cells.Item(25, 2).Value = "NEARProtocol"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(25, 4).Value = "'6.30"
$ws.Cells.Item(25, 5).Value = "  +12.04%  "
$ws.Cells.Item(26, 2).Value = "Litecoin"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(26, 4).Value = "'101.91"
$ws.Cells.Item(26, 5).Value = "  +15.01%  "
$ws.Cells.Item(27, 2).Value = "Aptos"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(27, 4).Value = "'12.71"
$ws.Cells.Item(27, 5).Value = "  +6.39%  "
$ws.Cells.Item(28, 2).Value = "WrappedeETH"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(28, 4).Value = "3.564.72"
$ws.Cells.Item(28, 5).Value = "  +0.18%  "
$ws.Cells.Item(29, 2).Value = "Hedera"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(29, 4).Value = "'0.150"
$ws.Cells.Item(29, 5).Value = "  +15.02%  "
$ws.Cells.Item(30, 2).Value = "Dai"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(30, 4).Value = "'0.999"
$ws.Cells.Item(30, 5).Value = "  +0.04%  "
$ws.Cells.Item(31, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(31, 4).Value = "'11.16"
$ws.Cells.Item(31, 5).Value = "  +15.73%  "
$ws.Cells.Item(32, 2).Value = "Cronos"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(32, 4).Value = "'0.190"
$ws.Cells.Item(32, 5).Value = "  -1.04%  "
$ws.Cells.Item(33, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(33, 4).Value = "'1.00"
$ws.Cells.Item(33, 5).Value = "  +0.47%  "
$ws.Cells.Item(34, 2).Value = "EthereumClassic"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(34, 4).Value = "'29.89"
$ws.Cells.Item(34, 5).Value = "  +6.27%  "
$ws.Cells.Item(35, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(35, 4).Value = "'0.545"
$ws.Cells.Item(35, 5).Value = "  +18.41%  "
$ws.Cells.Item(36, 2).Value = "RenderToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Cells.Item(36, 4).Value = "'7.95"
$ws.Cells.Item(36, 5).Value = "  +8.62%  "
$ws.Cells.Item(37, 2).Value = "PancakeSwap"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(37, 4).Value = "'2.14"
$ws.Cells.Item(37, 5).Value = "  +8.89%  "
$ws.Cells.Item(38, 2).Value = "Kaspa"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(38, 4).Value = "'0.161"
$ws.Cells.Item(38, 5).Value = "  +7.42%  "
$ws.Cells.Item(39, 2).Value = "Bittensor"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(39, 4).Value = "'527.40"
$ws.Cells.Item(39, 5).Value = "  +5.54%  "
$ws.Cells.Item(40, 2).Value = "VeChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(40, 4).Value = "'0.0455"
$ws.Cells.Item(40, 5).Value = "  +39.44%  "
$ws.Cells.Item(41, 2).Value = "Fetch.AI"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(41, 4).Value = "'1.34"
$ws.Cells.Item(41, 5).Value = "  +5.93%  "
$ws.Cells.Item(42, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(42, 4).Value = "'24.69"
$ws.Cells.Item(42, 5).Value = "  -0.88%  "
$ws.Cells.Item(43, 2).Value = "MantraDAO"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Cells.Item(43, 4).Value = "'3.80"
$ws.Cells.Item(43, 5).Value = "  +0.53%  "
$ws.Cells.Item(44, 2).Value = "ARBITRUM"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(44, 4).Value = "'0.841"
$ws.Cells.Item(44, 5).Value = "  +6.54%  "
$ws.Cells.Item(45, 2).Value = "dogwifhat"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(45, 4).Value = "'3.36"
$ws.Cells.Item(45, 5).Value = "  +2.92%  "
$ws.Cells.Item(46, 2).Value = "USDe"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(46, 4).Value = "'1.00"
$ws.Cells.Item(46, 5).Value = "  +0.04%  "
$ws.Cells.Item(47, 2).Value = "Stacks"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(47, 4).Value = "'2.07"
$ws.Cells.Item(47, 5).Value = "  +6.80%  "
$ws.Cells.Item(48, 2).Value = "Cosmos"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(48, 4).Value = "'7.93"
$ws.Cells.Item(48, 5).Value = "  +19.98%  "
$ws.Cells.Item(49, 2).Value = "Filecoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(49, 4).Value = "'5.14"
$ws.Cells.Item(49, 5).Value = "  +11.18%  "
$ws.Cells.Item(50, 2).Value = "OKB"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(50, 4).Value = "'50.96"
$ws.Cells.Item(50, 5).Value = "  +10.72%  "
$ws.Cells.Item(51, 2).Value = "ImmutableX"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(51, 4).Value = "'1.53"
$ws.Cells.Item(51, 5).Value = "  +11.61%  "
